$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the literal text representation of price/volume values (avoid Excel
# auto-converting numeric-looking strings and dropping trailing zeros, switching
# to scientific notation, etc.) by forcing Text format before assigning .Value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.381.64"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.568.60"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.21"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3758"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.18"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3399"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07576"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.137"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.99"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.973"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.926"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.566.61"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001129"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.87"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06754"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.57"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.199"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.385.07"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.711"
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.15"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.90"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.028"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.61"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.739.26"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.062"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9868"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.07"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.413"
$ws.Range("E36").Value = "  +9.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08463"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02490"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2287"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06480"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.412"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6306"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.21"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.98"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.802"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5923"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.073"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.263"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.77"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07327"
$ws.Range("E51").Value = "  +0.76%  "
